# Apply cryptos list price/volume refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to retain a literal text value (these Price cells are
    # plain General-format strings like "67.025.30" or "1.00" that Excel would
    # otherwise auto-coerce into a number, dropping formatting / precision).
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextCell "D2" "67.025.30"
$ws.Range("E2").Value = "  +0.13%  "

Set-TextCell "D3" "3.116.26"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").Value = "  -0.01%  "

Set-TextCell "D5" "580.44"
$ws.Range("E5").Value = "  -0.07%  "

Set-TextCell "D6" "173.48"
$ws.Range("E6").Value = "  +0.19%  "

Set-TextCell "D9" "6.41"
$ws.Range("E9").Value = "  -1.06%  "

$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("E11").Value = "  -0.73%  "

Set-TextCell "D12" "0.0000248"
$ws.Range("E12").Value = "  -0.50%  "

Set-TextCell "D13" "37.18"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("E14").Value = "  -1.55%  "

Set-TextCell "D15" "3.635.86"
$ws.Range("E15").Value = "  +0.44%  "

Set-TextCell "D16" "67.007.15"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("E17").Value = "  -0.94%  "

Set-TextCell "D18" "3.119.59"
$ws.Range("E18").Value = "  +0.33%  "

Set-TextCell "D19" "16.37"
$ws.Range("E19").Value = "  +1.86%  "

Set-TextCell "D20" "491.16"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("E21").Value = "  +5.40%  "

$ws.Range("E22").Value = "  -1.22%  "

Set-TextCell "D23" "84.11"
$ws.Range("E23").Value = "  +0.04%  "

Set-TextCell "D24" "13.21"
$ws.Range("E24").Value = "  +0.47%  "

Set-TextCell "D25" "2.28"
$ws.Range("E25").Value = "  -4.14%  "

Set-TextCell "D26" "10.38"
$ws.Range("E26").Value = "  +3.02%  "

Set-TextCell "D27" "1.00"
$ws.Range("E27").Value = "  +0.02%  "

Set-TextCell "D28" "7.90"
$ws.Range("E28").Value = "  -0.91%  "

$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("E30").Value = "  -0.53%  "

Set-TextCell "D31" "28.56"
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("E32").Value = "  -0.94%  "

Set-TextCell "D33" "0.0₃0946"
$ws.Range("E33").Value = "  -6.57%  "

Set-TextCell "D34" "0.999"
$ws.Range("E34").Value = "  -0.07%  "

Set-TextCell "D35" "5.87"
$ws.Range("E35").Value = "  -0.67%  "

$ws.Range("E36").Value = "  -2.21%  "

Set-TextCell "D37" "47.45"
$ws.Range("E37").Value = "  -1.30%  "

Set-TextCell "D38" "2.05"
$ws.Range("E38").Value = "  -3.17%  "

$ws.Range("E39").Value = "  -2.59%  "

Set-TextCell "D40" "0.123"
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("E41").Value = "  -1.61%  "

Set-TextCell "D42" "2.822.33"
$ws.Range("E42").Value = "  -0.42%  "

Set-TextCell "D43" "382.79"
$ws.Range("E43").Value = "  -0.48%  "

$ws.Range("E44").Value = "  -7.75%  "

Set-TextCell "D45" "0.0352"
$ws.Range("E45").Value = "  -2.69%  "

Set-TextCell "D46" "135.57"
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("E49").Value = "  -1.21%  "

$ws.Range("E50").Value = "  -0.85%  "

Set-TextCell "D51" "6.74"
$ws.Range("E51").Value = "  -0.97%  "

